$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Qualified Receiver Name" row label to plural "Qualified Receiver Names"
$ws.Range("A11").Value = "Qualified Receiver Names"

# Match the author's last selection state when they saved
$ws.Range("B22").Select() | Out-Null
